$d = $word.ActiveDocument

# The document ends with a trailing empty paragraph (<w:p/>). We replace it with two
# new paragraphs describing June 7th, 2022 (a diary-style entry), matching the other
# date/weather paragraph pairs already present in the document.

$lastIdx = $d.Paragraphs.Count
$pLast = $d.Paragraphs.Item($lastIdx)
$rLast = $pLast.Range
$insStart = $rLast.Start

# --- Paragraph A: the date line -----------------------------------------------------
# "2022年6月7日星期二"  (2022, June 7th, Tuesday)
$dateText = "2022年6月7日星期二"
$rLast.InsertAfter($dateText)

# --- Paragraph B: the weather / diary line ------------------------------------------
# Insert a paragraph break after the date line, then the weather/diary text.
$pA = $d.Paragraphs.Item($lastIdx)
$pA.Range.InsertParagraphAfter()
$pB = $d.Paragraphs.Item($lastIdx + 1)
$weatherText = "晴，今天是高考第一天，上午考语文，下午考数学，今天天气不错"
$pB.Range.InsertAfter($weatherText)

# --- Recreate the original run-splitting for paragraph A ----------------------------
# The source document types mixed Chinese/digit diary dates as several separate runs
# (digits vs. Chinese text groups are typed/committed separately). Recreate the same
# run boundaries by toggling a character property on/off across each tail sub-range,
# which forces Word to split runs at that boundary without altering the rendered text.
$paraAEnd = $pA.Range.End - 1
$dateBoundaries = @(1, 4, 5, 6, 7, 8, 11)
foreach ($off in $dateBoundaries) {
    $pos = $insStart + $off
    $rr = $d.Range($pos, $paraAEnd)
    $rr.Bold = 1
    $rr.Bold = 0
}

# --- Recreate the original run-splitting for paragraph B ----------------------------
$bStart = $pB.Range.Start
$paraBEnd = $pB.Range.End - 1
$weatherBoundaries = @(5)
foreach ($off in $weatherBoundaries) {
    $pos = $bStart + $off
    $rr = $d.Range($pos, $paraBEnd)
    $rr.Bold = 1
    $rr.Bold = 0
}
